# Color the USU-01 and USU-03 "Historias de Usuario" table rows green
# (RGB 92D050) across all three cells (Código, Historia de Usuario,
# Criterios de Aceptación), matching the author's "Ready" highlighting.

$d = $word.ActiveDocument

$green = 5296274  # wdColor BGR encoding of RGB(0x92, 0xD0, 0x50) -> w:color val="92D050"

$table = $d.Tables.Item(1)

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $codeCell = $table.Cell($r, 1)
    $codeText = $codeCell.Range.Text.TrimEnd([char]13, [char]7)
    if ($codeText -eq "USU-01" -or $codeText -eq "USU-03") {
        $table.Rows.Item($r).Range.Font.Color = $green
    }
}
